$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up the existing row 38 "Total" row: the blank placeholder cells
# (A38,B38,D38,E38,F38) are removed, leaving only C38 ("Total") and G38 (20000).
$ws.Range("A38").ClearContents()
$ws.Range("B38").ClearContents()
$ws.Range("D38").ClearContents()
$ws.Range("E38").ClearContents()
$ws.Range("F38").ClearContents()

# --- New order rows appended below the existing data (rows 39-56) ---

$rows = @(
    @{ R=39; A="2025-03-01 02:02:15"; B="camilo";  C="Sencilla";    D="Salchipapas";   E=2; F=10000; G=20000 },
    @{ R=40; Total=$true;                           C="Total";                                          G=20000 },

    @{ R=41; A="2025-03-01 02:08:43"; B="camilo";  C="Sencilla";    D="Salchipapas";   E=2; F=10000; G=20000 },
    @{ R=42; Total=$true;                           C="Total";                                          G=20000 },

    @{ R=43; A="2025-03-01 02:09:45"; B="camilo";  C="Sencilla";    D="Salchipapas";   E=2; F=10000; G=20000 },
    @{ R=44; Total=$true;                           C="Total";                                          G=20000 },

    @{ R=45; A="2025-03-01 02:19:10"; B="benichi"; C="Combinado";   D="Perro Caliente"; E=1; F=9000;  G=9000 },
    @{ R=46; A="2025-03-01 02:19:10"; B="benichi"; C="Gemelo";      D="Perro Caliente"; E=2; F=9000;  G=18000 },
    @{ R=47; A="2025-03-01 02:19:10"; B="benichi"; C="Suizo";       D="Perro Caliente"; E=3; F=12000; G=36000 },
    @{ R=48; Total=$true;                           C="Total";                                          G=63000 },

    @{ R=49; A="2025-03-01 02:25:40"; B="benichi"; C="Sencilla";    D="Salchipapas";   E=3; F=10000; G=30000 },
    @{ R=50; A="2025-03-01 02:25:40"; B="benichi"; C="Mixta";       D="Salchipapas";   E=3; F=17000; G=51000 },
    @{ R=51; A="2025-03-01 02:25:40"; B="benichi"; C="Combinada";   D="Salchipapas";   E=2; F=12000; G=24000 },
    @{ R=52; Total=$true;                           C="Total";                                          G=105000 },

    @{ R=53; A="2025-03-01 02:29:28"; B="benichi"; C="Ranchera";    D="Salchipapas";   E=3; F=15000; G=45000 },
    @{ R=54; A="2025-03-01 02:29:28"; B="benichi"; C="Combinada";   D="Salchipapas";   E=2; F=12000; G=24000 },
    @{ R=55; A="2025-03-01 02:29:28"; B="benichi"; C="Salchipollo"; D="Salchipapas";   E=2; F=15000; G=30000 },
    @{ R=56; Total=$true; Blank=$true;              C="Total";                                          G=99000 }
)

foreach ($row in $rows) {
    $r = $row.R

    if ($row.ContainsKey("A")) {
        $ws.Cells.Item($r, 1).Value = $row.A
    } elseif ($row.Blank) {
        $ws.Cells.Item($r, 1).ClearContents()
    }

    if ($row.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $row.B
    } elseif ($row.Blank) {
        $ws.Cells.Item($r, 2).ClearContents()
    }

    $ws.Cells.Item($r, 3).Value = $row.C

    if ($row.ContainsKey("D")) {
        $ws.Cells.Item($r, 4).Value = $row.D
    } elseif ($row.Blank) {
        $ws.Cells.Item($r, 4).ClearContents()
    }

    if ($row.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $row.E
    } elseif ($row.Blank) {
        $ws.Cells.Item($r, 5).ClearContents()
    }

    if ($row.ContainsKey("F")) {
        $ws.Cells.Item($r, 6).Value = $row.F
    } elseif ($row.Blank) {
        $ws.Cells.Item($r, 6).ClearContents()
    }

    $ws.Cells.Item($r, 7).Value = $row.G
}
